$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (B:TB, C:d2S, D:K, E:IP, F:Win, G:sum)
# F (Win) is unchanged; G (sum) = B + C + D + E
$data = @{
    2  = @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 3.009163075608874)
    3  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    4  = @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 11.54378914222666)
    5  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    6  = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 27.82738278199502)
    7  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    8  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    9  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    10 = @(0.6606524410359556, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 2.214453472130288)
    11 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    12 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    13 = @(1.455362044514542, 117.745847958593, 0.7527432677738641, 10.19245300693656, 130.146406277818)
    14 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    15 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302)
    16 = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478)
    17 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    18 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
